$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Change 1: "АСТАНА 202" + "3" (two runs) -> single run "АСТАНА 2023"
# ---------------------------------------------------------------------------
$rngYear = $d.Content
$rngYear.Find.Execute("АСТАНА 2023", $false, $false, $false, $false, $false, `
                       $true, 1, $false, "АСТАНА 2023", 2) | Out-Null

# ---------------------------------------------------------------------------
# Change 2: remove the empty paragraph that sits between the closing "."
# paragraph and the "1. Общие положения" heading.
# ---------------------------------------------------------------------------
$rngHeading1 = $d.Content
$rngHeading1.Find.Execute("1. Общие положения", $false, $false, $false, $false, `
                           $false, $true, 1, $false, "", 0) | Out-Null
$heading1 = $rngHeading1.Paragraphs(1)
$emptyBefore1 = $heading1.Previous()
$emptyBefore1.Range.Delete()

# ---------------------------------------------------------------------------
# Change 3: insert one more empty paragraph (cloned from the preceding
# empty paragraph) right before the "2. История появления" heading.
# ---------------------------------------------------------------------------
$rngHeading2 = $d.Content
$rngHeading2.Find.Execute("2. История появления", $false, $false, $false, `
                           $false, $false, $true, 1, $false, "", 0) | Out-Null
$heading2 = $rngHeading2.Paragraphs(1)
$lastEmpty2 = $heading2.Previous()
$lastEmpty2.Range.InsertParagraphAfter()

# ---------------------------------------------------------------------------
# Change 4: give the (still empty) run right after the "2. История
# появления" heading the same direct formatting as the heading's mark.
# ---------------------------------------------------------------------------
$rngHeading2b = $d.Content
$rngHeading2b.Find.Execute("2. История появления", $false, $false, $false, `
                            $false, $false, $true, 1, $false, "", 0) | Out-Null
$heading2b = $rngHeading2b.Paragraphs(1)
$afterHeading2 = $heading2b.Next()
$afterRng = $afterHeading2.Range
$afterRng.Font.NameAscii = "Times New Roman"
$afterRng.Font.Name = "Times New Roman"
$afterRng.Bold = $true
$afterRng.BoldBi = $true
$afterRng.Font.Size = 14
$afterRng.Font.SizeBi = 14
$afterRng.LanguageID = "ru-RU"
